$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Setting")
Write-Output "noop"
